$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# Remove the EmailAddress / ConfirmEmailAddress hyperlinks, then delete those
# two columns entirely so Password / ConfirmPassword shift left into I:J.
$ws.Hyperlinks.Delete()
$ws.Range("I1:J1").EntireColumn.Delete()

# The "Hyperlink" cell style is no longer used anywhere in the sheet - drop it.
$wb.Styles("Hyperlink").Delete()

$ws.Range("G18").Select()
